$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$noteButtercup = "John was not feeling well.  He had  diarrhea.  This sentence contains the word buttercup exactly one time."
$noteRemoved   = "John was not feeling well.  He had  diarrhea.  This sentence contains the word  exactly one time."

$gValues = @(15, 16, 17, 18, 19, 20, 21, 15)

for ($r = 2; $r -le 9; $r++) {
    if ($r -eq 5) {
        $ws.Cells.Item($r, 1).Value = $noteRemoved
    } else {
        $ws.Cells.Item($r, 1).Value = $noteButtercup
    }

    $ws.Cells.Item($r, 2).Value = "Stevens, Rick"
    $ws.Cells.Item($r, 3).Value = $null
    $ws.Cells.Item($r, 4).Value = "Harvest"
    $ws.Cells.Item($r, 5).Value = 0.35416666666666669
    $ws.Cells.Item($r, 6).Value = 0.60416666666666663
    $ws.Cells.Item($r, 7).Value = $gValues[$r - 2]
    $ws.Cells.Item($r, 8).Value = "Staff Jones"
}

$ws.Range("A5").Select()
